$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header cells: A1 becomes "maturities", B1 becomes "quotes"
$ws.Range("A1").Value = "maturities"
$ws.Range("B1").Value = "quotes"

# Convert the month-count values in A2:A32 into "<n>M" text labels
for ($r = 2; $r -le 32; $r++) {
    $n = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = "$([int]$n)M"
}

# Update the selection to E13 to match the workbook's last-saved selection
$ws.Range("E13").Select()
